# anggota_keluar.docx — T1706/T1705/T1704
# "Fix akad keluar anggota, jumlah di menu simpanan, jumlah di menu tabungan"
#
# Content changes applied:
#   1. "SIMPANAN SUKARELA" -> "SIMPANAN KHUSUS" (only the SUKARELA -> KHUSUS
#      portion is rewritten, keeping it as a distinct trailing run so the
#      paragraph ends up as two runs: "SIMPANAN " + "KHUSUS").
#   2. Remove the stray "_GoBack" bookmark that wrapped the "TOTAL" paragraph.

$d = $word.ActiveDocument

# --- 1. SIMPANAN SUKARELA -> SIMPANAN KHUSUS -----------------------------
# Match case-sensitively so we only touch the heading text "SUKARELA" and
# leave the ${jumlah_simpanan_sukarela} merge field placeholder untouched.
$d.Content.Find.Execute("SUKARELA", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "KHUSUS", 2)

# Re-apply (no-op) direct formatting to the freshly inserted "KHUSUS" text so
# it is kept as its own run instead of being coalesced back into the
# preceding "SIMPANAN " run.
$khusus = $d.Content
$khusus.Find.Execute("KHUSUS")
$khusus.Bold = 1
$khusus.Bold = 0

# --- 2. Drop the _GoBack bookmark around the TOTAL paragraph -------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()
